$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("general")
$ws.Rows.Item(6).Insert()
$ws.Rows.Item(6).Insert()
$ws.Rows.Item(1048577).EntireRow.Delete()
"done" | Out-Host
